$wb = $excel.ActiveWorkbook

# Keep calculation manual so the existing cached formula results (e.g. B6)
# are left untouched by this edit, matching the source change which only
# touched the raw cell values/types.
$excel.Calculation = -4135

$ws = $wb.Worksheets.Item(1)

# Row 6 currently has numeric values 436 (D6) and 1705 (E6).
# Replace them with the text values "4:36" and "17:05" respectively,
# explicitly formatted (Calibri 11, black) as rich-text runs so the
# formatting travels with the shared string rather than the cell style.

$ws.Range("D6").Value = "4:36"
$d6a = $ws.Range("D6").Characters(1, 2)
$d6a.Font.Name = "Calibri"
$d6a.Font.Size = 11
$d6a.Font.Color = 0
$d6b = $ws.Range("D6").Characters(3, 2)
$d6b.Font.Name = "Calibri"
$d6b.Font.Size = 11
$d6b.Font.Color = 0

$ws.Range("E6").Value = "17:05"
$e6a = $ws.Range("E6").Characters(1, 3)
$e6a.Font.Name = "Calibri"
$e6a.Font.Size = 11
$e6a.Font.Color = 0
$e6b = $ws.Range("E6").Characters(4, 2)
$e6b.Font.Name = "Calibri"
$e6b.Font.Size = 11
$e6b.Font.Color = 0
